# reworked service to map to Stats, calculate, save, and then map to Roto
# and also fixed excel writer to have space

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fixed the "__" placeholder column header to be a single space.
$ws.Range("F1").Value = " "

# Recalculated stats: rankings / totals / hitting / pitching shuffled
# around after recompute. Rows 2-5's team names (rank 1-4) kept their
# place; rows 6-15 got reshuffled along with the recalculated numbers.

$numeric = @(
    @{ Row = 2;  Rank = 1.0;  Total = 133.0; Hitting = 60.0; Pitching = 73.0 },
    @{ Row = 3;  Rank = 2.0;  Total = 124.0; Hitting = 58.0; Pitching = 66.0 },
    @{ Row = 4;  Rank = 3.0;  Total = 114.0; Hitting = 56.0; Pitching = 58.0 },
    @{ Row = 5;  Rank = 4.0;  Total = 99.0;  Hitting = 53.0; Pitching = 46.0 },
    @{ Row = 6;  Rank = 5.0;  Total = 97.0;  Hitting = 46.0; Pitching = 51.0 },
    @{ Row = 7;  Rank = 6.0;  Total = 89.0;  Hitting = 33.0; Pitching = 56.0 },
    @{ Row = 8;  Rank = 7.0;  Total = 88.0;  Hitting = 34.0; Pitching = 54.0 },
    @{ Row = 9;  Rank = 8.0;  Total = 83.0;  Hitting = 36.0; Pitching = 47.0 },
    @{ Row = 10; Rank = 9.0;  Total = 78.0;  Hitting = 44.0; Pitching = 34.0 },
    @{ Row = 11; Rank = 10.5; Total = 72.0;  Hitting = 47.0; Pitching = 25.0 },
    @{ Row = 12; Rank = 10.5; Total = 72.0;  Hitting = 37.0; Pitching = 35.0 },
    @{ Row = 13; Rank = 12.0; Total = 67.0;  Hitting = 36.0; Pitching = 31.0 },
    @{ Row = 14; Rank = 13.0; Total = 64.0;  Hitting = 43.0; Pitching = 21.0 },
    @{ Row = 15; Rank = 14.0; Total = 62.0;  Hitting = 40.0; Pitching = 22.0 }
)

foreach ($r in $numeric) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Rank
    $ws.Cells.Item($row, 3).Value = $r.Total
    $ws.Cells.Item($row, 4).Value = $r.Hitting
    $ws.Cells.Item($row, 5).Value = $r.Pitching
}

# Team order (column B) for the rows whose rank moved. Rows 2-5 (and
# their names) are untouched so the embedded special characters stay
# exactly as authored.
$names = @(
    @{ Row = 6;  Name = "GOD WILLS IT" },
    @{ Row = 7;  Name = "Swampnuts" },
    @{ Row = 8;  Name = "Splitfinger Skadoosh" },
    @{ Row = 9;  Name = "Epic7" },
    @{ Row = 10; Name = "MillerTime" },
    @{ Row = 11; Name = "confusion" },
    @{ Row = 12; Name = "SmokeWalkers" },
    @{ Row = 13; Name = "Mac" },
    @{ Row = 14; Name = "Corbin Copy" },
    @{ Row = 15; Name = "DJ's Quality Team" }
)

foreach ($n in $names) {
    $ws.Cells.Item($n.Row, 2).Value = $n.Name
}
